$d = $word.ActiveDocument

# Locate the "Trial here" run and replace its text with "Game".
$r = $d.Content
$found = $r.Find.Execute("Trial here", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    $start = $r.Start
    $r.Text = "Game"

    # Force the engine to emit "Game" as its own run (matching the target
    # markup, which splits "| Trial here: " into "| " / "Game" / ": ")
    # by explicitly (re)stamping run-level character formatting on just
    # the inserted word.
    $gameRange = $d.Range($start, $start + 4)
    $gameRange.Font.Bold = $true
    $gameRange.Font.Bold = $false
}
